# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scrape output, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Row => new F-column value, for worksheet "展览"
$exhibitionUpdates = @{
    2  = 258
    3  = 1095
    5  = 448
    6  = 85
    7  = 576
    8  = 81
    9  = 6882
    10 = 166
    15 = 1117
    16 = 16338
    17 = 3
    20 = 339
    21 = 192
    23 = 11448
    25 = 1079
    26 = 4505
    27 = 360
    30 = 849
    31 = 324
    32 = 144
}

# Row => new F-column value, for worksheet "全部类型"
$allTypesUpdates = @{
    2  = 258
    3  = 1095
    5  = 448
    6  = 85
    7  = 576
    9  = 81
    10 = 6882
    11 = 166
    17 = 1117
    18 = 16338
    19 = 3
    22 = 339
    23 = 192
    27 = 11448
    29 = 1079
    30 = 4505
    31 = 360
    34 = 849
    35 = 324
    36 = 144
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
